# Initial processing of 2023 data
# Remove the "Not available via API" comments from rows 2 and 4 (column E)
# since that data has now been obtained / processed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("observed_stns")

$ws.Range("E2").ClearContents()
$ws.Range("E4").ClearContents()

# Update the active selection to reflect where the editor's cursor ended up
$ws.Range("E6").Select()
